# Edit script for tut06/output/2001EE28.xlsx
# - Change date strings in column A from DD/MM/YYYY to DD-MM-YYYY for rows 3..21
# - Update numeric values in columns D, E, G, H for the specific rows that changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (slash -> dash)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "'" + $dates[$row]
}

# Updated numeric values for columns D (4), E (5), G (7), H (8)
# Row 3: D 0->1, G 0->1 (E,H unchanged)
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1

# Rows 4,5,6,11,16: D 0->1, E 0->1, H 1->0
$rowsSet1 = @(4, 5, 6, 11, 16)
foreach ($row in $rowsSet1) {
    $ws.Cells.Item($row, 4).Value = 1
    $ws.Cells.Item($row, 5).Value = 1
    $ws.Cells.Item($row, 8).Value = 0
}
